# Atualização de bases das ligas, do dia: 27-04-2024 às 11:27
#
# The "Romania Liga I" fixtures sheet got refreshed match odds data.
# Rows 235 & 237 (match ids 6861095 / 6870268) and rows 238 & 239
# (match ids 6852370 / 6836277) had their entire data records (columns
# B:AB) swapped with one another, and a handful of odds cells in rows
# 292, 293 and 296 were simply updated in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $row1, $row2, $startCol, $endCol) {
    $rng1 = $ws.Range($ws.Cells.Item($row1, $startCol), $ws.Cells.Item($row1, $endCol))
    $rng2 = $ws.Range($ws.Cells.Item($row2, $startCol), $ws.Cells.Item($row2, $endCol))
    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2
    $rng1.Value = $vals2
    $rng2.Value = $vals1
}

# Columns B (id) through AB (PL_AhUnder) = columns 2..28
$firstCol = 2
$lastCol = 28

# Swap full records between row 235 <-> row 237
Swap-Rows $ws 235 237 $firstCol $lastCol

# Swap full records between row 238 <-> row 239
Swap-Rows $ws 238 239 $firstCol $lastCol

# Row 292: in-place odds refresh
$ws.Range("Q292").Value = 2
$ws.Range("R292").Value = 1.85
$ws.Range("S292").Value = 2
$ws.Range("T292").Value = 1.825
$ws.Range("U292").Value = 2.025

# Row 293: in-place odds refresh
$ws.Range("Q293").Value = 1.875
$ws.Range("R293").Value = 1.975
$ws.Range("T293").Value = 1.8
$ws.Range("U293").Value = 2.05

# Row 296: in-place odds refresh
$ws.Range("N296").Value = 4.2
$ws.Range("T296").Value = 1.85
$ws.Range("U296").Value = 2
